# Refresh the cryptocurrency snapshot table (Coin / Link / Price / Volume(1h))
# with the latest pulled values, matching the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" cells are numeric-looking strings that must stay literal
# text (e.g. trailing zeros like "1.00" or "0.0000270"). Force text format on
# just those cells so Excel does not silently coerce them to numbers.
$textForcedCells = @("D13", "D25", "D27", "D29", "D43", "D48")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '66.835.71'
$ws.Range("E2").Value = '  +5.81%  '

# Row 3
$ws.Range("D3").Value = '3.538.03'
$ws.Range("E3").Value = '  +9.53%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '187.61'
$ws.Range("E5").Value = '  +9.69%  '

# Row 6
$ws.Range("D6").Value = '558.72'
$ws.Range("E6").Value = '  +5.85%  '

# Row 7
$ws.Range("D7").Value = '3.532.76'
$ws.Range("E7").Value = '  +9.40%  '

# Row 8
$ws.Range("D8").Value = '0.616'
$ws.Range("E8").Value = '  +4.07%  '

# Row 9
$ws.Range("E9").Value = '  -0.09%  '

# Row 10
$ws.Range("D10").Value = '0.633'
$ws.Range("E10").Value = '  +4.88%  '

# Row 11
$ws.Range("E11").Value = '  +13.60%  '

# Row 12
$ws.Range("D12").Value = '54.66'
$ws.Range("E12").Value = '  +3.01%  '

# Row 13
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  +6.84%  '

# Row 14
$ws.Range("D14").Value = '9.37'
$ws.Range("E14").Value = '  +2.90%  '

# Row 15
$ws.Range("D15").Value = '4.102.56'
$ws.Range("E15").Value = '  +9.42%  '

# Row 16
$ws.Range("D16").Value = '3.539.51'
$ws.Range("E16").Value = '  +9.55%  '

# Row 17
$ws.Range("E17").Value = '  +4.73%  '

# Row 18
$ws.Range("D18").Value = '66.876.40'

# Row 19
$ws.Range("D19").Value = '18.21'
$ws.Range("E19").Value = '  +6.19%  '

# Row 20
$ws.Range("D20").Value = '11.98'
$ws.Range("E20").Value = '  +8.68%  '

# Row 21
$ws.Range("D21").Value = '0.995'
$ws.Range("E21").Value = '  +3.10%  '

# Row 22
$ws.Range("D22").Value = '432.14'

# Row 23
$ws.Range("D23").Value = '4.13'
$ws.Range("E23").Value = '  +10.01%  '

# Row 24
$ws.Range("D24").Value = '85.06'
$ws.Range("E24").Value = '  +4.87%  '

# Row 25
$ws.Range("D25").Value = '4.10'
$ws.Range("E25").Value = '  +1.73%  '

# Row 26
$ws.Range("D26").Value = '11.07'
$ws.Range("E26").Value = '  -0.58%  '

# Row 27
$ws.Range("D27").Value = '2.90'
$ws.Range("E27").Value = '  +9.97%  '

# Row 28
$ws.Range("D28").Value = '6.13'
$ws.Range("E28").Value = '  -0.11%  '

# Row 29
$ws.Range("D29").Value = '12.20'
$ws.Range("E29").Value = '  +8.48%  '

# Row 30
$ws.Range("D30").Value = '9.14'
$ws.Range("E30").Value = '  +11.75%  '

# Row 31
$ws.Range("D31").Value = '30.31'
$ws.Range("E31").Value = '  +6.53%  '

# Row 32
$ws.Range("D32").Value = '645.53'
$ws.Range("E32").Value = '  +1.24%  '

# Row 33
$ws.Range("D33").Value = '6.59'
$ws.Range("E33").Value = '  +2.41%  '

# Row 34
$ws.Range("D34").Value = '11.74'
$ws.Range("E34").Value = '  +4.88%  '

# Row 35
$ws.Range("E35").Value = '  +6.05%  '

# Row 36
$ws.Range("D36").Value = '59.66'
$ws.Range("E36").Value = '  +4.85%  '

# Row 37
$ws.Range("D37").Value = '0.152'
$ws.Range("E37").Value = '  +23.91%  '

# Row 38
$ws.Range("E38").Value = '  +15.30%  '

# Row 39
$ws.Range("D39").Value = '38.61'
$ws.Range("E39").Value = '  +5.62%  '

# Row 40
$ws.Range("E40").Value = '  +0.06%  '

# Row 41
$ws.Range("E41").Value = '  +3.57%  '

# Row 42
$ws.Range("E42").Value = '  +14.44%  '

# Row 43
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.11%  '

# Row 44
$ws.Range("D44").Value = '3.051.79'
$ws.Range("E44").Value = '  +6.03%  '

# Row 45
$ws.Range("E45").Value = '  +4.03%  '

# Row 46
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.39'
$ws.Range("E46").Value = '  +10.74%  '

# Row 47
$ws.Range("B47").Value = 'ThetaToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D47").Value = '2.87'
$ws.Range("E47").Value = '  +11.89%  '

# Row 48
$ws.Range("D48").Value = '2.80'
$ws.Range("E48").Value = '  +4.47%  '

# Row 49
$ws.Range("D49").Value = '0.0419'
$ws.Range("E49").Value = '  +6.51%  '

# Row 50
$ws.Range("D50").Value = '0.132'
$ws.Range("E50").Value = '  +5.64%  '

# Row 51
$ws.Range("D51").Value = '8.65'
$ws.Range("E51").Value = '  +10.49%  '
